$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.034668445587158
$ws.Range("B1").Value = 1.018976926803589
$ws.Range("C1").Value = 0.7948615550994873
$ws.Range("D1").Value = 5.245019912719727
$ws.Range("E1").Value = 2.015298843383789
